$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark from the very start of the document
#    (paragraph 1). It reappears later, near the end of paragraph 6.
# ---------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------
# 2. Paragraph 3 ("The upper part displays ... over the course of
#    time."): drop everything from "Each segment ..." to the end of
#    the paragraph, leaving just the introductory sentence.
# ---------------------------------------------------------------------
$found = $d.Content.Find.Execute(
    "Each segment represents a Month, during which the consumption is evaluated. A red highlight is the result of a monthly consumption of more than 150% of the designated Quota. The colors orange, dark green and light green are each representative of 150% to 110%, 110% to 70% and below 70%.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ---------------------------------------------------------------------
# 3. Paragraph 4 (previously "The grey graph ..."): becomes a new
#    paragraph explaining the X/Y axes.
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$rng4 = $p4.Range
$rng4.End = $rng4.End - 1
$rng4.Text = "The X-axis represents the job" + [char]8217 + "s date, while the Y-axis indicates the accumulated CPUhours. Each segment represents a Month, during which the consumption is evaluated. Each bar is highlighted in a color corresponding to how much of the monthly quota was occupied."

# ---------------------------------------------------------------------
# 4. Paragraph 5 (previously empty): gets the old "grey graph" text
#    (minus the "higher lower"/"black graph" sentences) plus a new
#    closing sentence about the overall Efficiency.
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$rng5 = $p5.Range
$nl = [char]11
$rng5.Text = "The grey graph signifies the occupied corehours, to the same effect, the yellow one does for the used corehours." + $nl + "The Quotient between used to occupied cpuhours is the overall Efficiency."
# the paragraph started out empty, so the freshly created run has no
# rPr of its own yet; stamp the same language the rest of the
# document uses so it matches its sibling runs.
$d.Paragraphs(5).Range.LanguageID = "en-US"

# ---------------------------------------------------------------------
# 5. Paragraph 6 ("The Second graph ..."): reworded, and the
#    "_GoBack" bookmark is re-inserted right before the final period.
# ---------------------------------------------------------------------
$p6 = $d.Paragraphs(6)
$rng6 = $p6.Range
$rng6.End = $rng6.End - 1
$rng6.Text = "The Second displays a purple dot for the efficiency of every job and a red dot for the average Efficiency per day."

$p6again = $d.Paragraphs(6)
$searchRng6 = $p6again.Range
$found6 = $searchRng6.Find.Execute("day.")
$bmPos = $searchRng6.Start + 3
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

# ---------------------------------------------------------------------
# 6. Paragraph 7 ("Only days, where Jobs ..."): "where" -> "on which".
# ---------------------------------------------------------------------
$rng7 = $d.Paragraphs(7).Range
$found7 = $rng7.Find.Execute("where", $false, $false, $false, $false, $false, $true, 1, $false, "on which", 2)

Write-Host "Done"
